$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Swap the "Python" / "node.js" rows and update node.js/Python request counts.
$ws.Range("A3").Value = "node.js"
$ws.Range("B3").Value = 2076

$ws.Range("A4").Value = "Python"
$ws.Range("B4").Value = 2112

# Reproduce the final selection state: the whole of row 5 selected (as if the
# row header for row 5 was clicked), with A5 as the active cell.
$ws.Rows("5").Select()
